$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store plain text exactly as given, regardless of
# whether the text "looks" numeric. We prefix with a text-marker apostrophe
# (mirrors typing '123 in Excel) and then reset the resulting style back to
# the default "Normal" so we do not leave a stray quote-prefix/format style
# on the cell (only its shared-string content should change, per the diff).
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

# U+2083 SUBSCRIPT THREE used in the PEPE price (row 38). Built via string
# interpolation (not the "+" operator) because this engine's "+" silently
# does numeric addition when both operands look like numbers/digits.
$sub3 = [char]0x2083

Set-TextCell "D2" '61.659.36'
Set-TextCell "E2" '  -1.57%  '
Set-TextCell "D3" '3.041.02'
Set-TextCell "E3" '  -4.64%  '
Set-TextCell "E4" '  +0.06%  '
Set-TextCell "D5" '580.56'
Set-TextCell "E5" '  -1.54%  '
Set-TextCell "D6" '130.08'
Set-TextCell "E6" '  -4.63%  '
Set-TextCell "E7" '  +0.03%  '
Set-TextCell "D8" '3.043.14'
Set-TextCell "E8" '  -4.44%  '
Set-TextCell "D9" '0.501'
Set-TextCell "E9" '  -1.36%  '
Set-TextCell "D10" '0.138'
Set-TextCell "E10" '  -3.20%  '
Set-TextCell "D11" '5.23'
Set-TextCell "E11" '  -1.02%  '
Set-TextCell "D12" '0.438'
Set-TextCell "E12" '  -3.58%  '
Set-TextCell "D13" '0.0000232'
Set-TextCell "E13" '  -1.71%  '
Set-TextCell "D14" '33.36'
Set-TextCell "E14" '  -0.17%  '
Set-TextCell "E15" '  +0.81%  '
Set-TextCell "D16" '3.550.94'
Set-TextCell "E16" '  -4.39%  '
Set-TextCell "D17" '61.754.70'
Set-TextCell "E17" '  -1.42%  '
Set-TextCell "D18" '3.045.35'
Set-TextCell "E18" '  -4.59%  '
Set-TextCell "D19" '6.36'
Set-TextCell "E19" '  -2.74%  '
Set-TextCell "D20" '447.22'
Set-TextCell "E20" '  -2.19%  '
Set-TextCell "D21" '13.46'
Set-TextCell "E21" '  -3.51%  '
Set-TextCell "D22" '0.669'
Set-TextCell "E22" '  -5.06%  '
Set-TextCell "D23" '7.32'
Set-TextCell "E23" '  -4.19%  '
Set-TextCell "D24" '80.82'
Set-TextCell "E24" '  -3.35%  '
Set-TextCell "D25" '12.78'
Set-TextCell "E25" '  -3.74%  '
Set-TextCell "E26" '  +0.06%  '
Set-TextCell "D27" '0.999'
Set-TextCell "E27" '  -0.13%  '
Set-TextCell "E28" '  -5.08%  '
Set-TextCell "D29" '2.00'
Set-TextCell "E29" '  -1.25%  '
Set-TextCell "D30" '7.40'
Set-TextCell "E30" '  -4.97%  '
Set-TextCell "D31" '6.46'
Set-TextCell "E31" '  -5.93%  '
Set-TextCell "D32" '25.85'
Set-TextCell "E32" '  -5.56%  '
Set-TextCell "D33" '0.0968'
Set-TextCell "E33" '  -6.38%  '
Set-TextCell "E34" '  -2.81%  '
Set-TextCell "D35" '0.972'
Set-TextCell "E35" '  -6.19%  '
Set-TextCell "D36" '5.69'
Set-TextCell "E36" '  -3.74%  '
Set-TextCell "D37" '50.28'
Set-TextCell "E37" '  -1.67%  '
Set-TextCell "D38" "0.0${sub3}0699"
Set-TextCell "E38" '  -0.47%  '
Set-TextCell "E39" '  -3.56%  '
Set-TextCell "D40" '7.90'
Set-TextCell "E40" '  -1.42%  '
Set-TextCell "E41" '  -2.54%  '
Set-TextCell "B42" 'dogwifhat'
Set-TextCell "C42" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell "D42" '2.52'
Set-TextCell "E42" '  -7.59%  '
Set-TextCell "B43" 'Bittensor'
Set-TextCell "C43" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell "D43" '376.68'
Set-TextCell "E43" '  -5.91%  '
Set-TextCell "D44" '2.685.20'
Set-TextCell "E44" '  -5.41%  '
Set-TextCell "D45" '0.998'
Set-TextCell "E45" '  -0.03%  '
Set-TextCell "D46" '123.12'
Set-TextCell "E46" '  -1.57%  '
Set-TextCell "D47" '0.238'
Set-TextCell "E47" '  -4.80%  '
Set-TextCell "D48" '34.13'
Set-TextCell "E48" '  -6.50%  '
Set-TextCell "D49" '2.00'
Set-TextCell "E49" '  -6.50%  '
Set-TextCell "E50" '  -3.03%  '
Set-TextCell "D51" '23.80'
Set-TextCell "E51" '  -6.95%  '
